$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content edits -------------------------------------------------------
# B3: "송여령" -> "O"
$ws.Range("B3").Value = "O"

# B4: "송여령`n(쿠키런B)" -> "O"  (the now-unused "송여령(쿠키런B)" shared
# string is dropped automatically once nothing references it any more)
$ws.Range("B4").Value = "O"

# B14: "송여령" -> "송여령 ing"
$ws.Range("B14").Value = "송여령 ing"

# B15: empty -> new note "O? 렉 많음"
$ws.Range("B15").Value = "O? 렉 많음"

# --- Selection / scroll position ------------------------------------------
# Selection moves from B14 to C15, and the window is scrolled back so A1 is
# the top-left visible cell again (was topLeftCell="A10").
$ws.Range("A1").Select()
$ws.Range("C15").Select()

# --- Row height touch-ups ---------------------------------------------------
# Row 4 no longer needs a custom (wrapped-text) height once it only holds "O".
$ws.Rows.Item(4).AutoFit()

# The remaining rows keep a custom height, just a touch taller than before
# (consistent with the refreshed default row height/font metrics).
$rowHeights = @{
    2  = 68
    5  = 34
    6  = 51
    7  = 34
    8  = 34
    9  = 34
    11 = 51
    12 = 34
    18 = 51
    19 = 51
    21 = 51
    24 = 51
    25 = 68
    26 = 51
    27 = 34
    29 = 34
    30 = 51
    31 = 68
}
foreach ($r in $rowHeights.Keys) {
    $ws.Rows.Item($r).RowHeight = $rowHeights[$r]
}
